$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Inventory" to "Sheet1"
$ws.Name = "Sheet1"

# Style the header row (A1:G1): bold font, thin box border, centered
# horizontally and top-aligned vertically.
$hdr = $ws.Range("A1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Copy the computed header format onto the rest of the header row so every
# cell shares a single style entry instead of each creating its own.
$hdr.Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Drop the duplicated "Delivery Date / Invoice Total" rows (old rows 6-7
# are removed first so the earlier rows 3-4 keep their original numbers).
$ws.Rows("6:7").Delete()
$ws.Rows("3:4").Delete()

# Consolidate the extracted label/value pairs into a single row under the
# header, combining label and value into one text string per cell.
$ws.Range("B2").Value = "Delivery Date: 2024-10-12"
$ws.Range("E2").Value = "Invoice Total: 500"
